# Updates "Price" (column D) and "Volume(1h)" (column E) figures in the
# cryptos worksheet to the latest scraped values.
#
# A handful of the new Price values are plain decimal numbers (e.g. "1.00",
# "318.71", "0.0869"). If such a string is assigned to a Range.Value as-is,
# Excel happily reinterprets it as a genuine number (losing trailing zeros
# and introducing floating point noise), which is not what we want here -
# the sheet stores these figures as plain text. So for any new value that
# would parse as a number we force the cell to Text format first via
# NumberFormat = "@" before writing the string, which keeps it as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    if ($text -match '^[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $text
}

Set-TextValue "D2"  "41.653.65"

Set-TextValue "D3"  "2.473.67"
Set-TextValue "E3"  "  -0.18%  "

Set-TextValue "E4"  "  +0.08%  "

Set-TextValue "D5"  "318.71"
Set-TextValue "E5"  "  +1.78%  "

Set-TextValue "D6"  "92.91"
Set-TextValue "E6"  "  -0.69%  "

Set-TextValue "E7"  "  +1.93%  "

Set-TextValue "E8"  "  +0.04%  "

Set-TextValue "E9"  "  +2.39%  "

Set-TextValue "D10" "0.0869"
Set-TextValue "E10" "  +10.75%  "

Set-TextValue "E11" "  +0.81%  "

Set-TextValue "D13" "2.854.76"
Set-TextValue "E13" "  -0.16%  "

Set-TextValue "D14" "6.91"
Set-TextValue "E14" "  +1.31%  "

Set-TextValue "E15" "  -2.32%  "

Set-TextValue "D16" "2.470.33"
Set-TextValue "E16" "  +0.50%  "

Set-TextValue "D17" "0.789"
Set-TextValue "E17" "  +3.42%  "

Set-TextValue "D18" "41.620.90"
Set-TextValue "E18" "  +0.24%  "

Set-TextValue "D19" "0.0₃0958"
Set-TextValue "E19" "  +2.30%  "

Set-TextValue "D20" "6.49"
Set-TextValue "E20" "  +1.55%  "

Set-TextValue "D21" "71.34"
Set-TextValue "E21" "  -0.31%  "

Set-TextValue "D22" "11.48"
Set-TextValue "E22" "  +1.59%  "

Set-TextValue "D23" "241.30"
Set-TextValue "E23" "  +2.01%  "

Set-TextValue "E24" "  +1.56%  "

Set-TextValue "E25" "  +1.86%  "

Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  -0.02%  "

Set-TextValue "D27" "24.85"
Set-TextValue "E27" "  +0.07%  "

Set-TextValue "D28" "2.28"
Set-TextValue "E28" "  +3.74%  "

Set-TextValue "E29" "  +2.75%  "

Set-TextValue "D30" "36.67"
Set-TextValue "E30" "  +2.41%  "

Set-TextValue "D31" "157.21"
Set-TextValue "E31" "  -0.70%  "

Set-TextValue "E32" "  +0.97%  "

Set-TextValue "E33" "  +0.03%  "

Set-TextValue "E34" "  +2.36%  "

Set-TextValue "E35" "  +0.08%  "

Set-TextValue "E36" "  -0.17%  "

Set-TextValue "E37" "  +0.07%  "

Set-TextValue "E38" "  +0.60%  "

Set-TextValue "E39" "  +1.51%  "

Set-TextValue "E40" "  -1.77%  "

Set-TextValue "D41" "4.02"
Set-TextValue "E41" "  -2.57%  "

Set-TextValue "E42" "  +2.04%  "

Set-TextValue "D43" "1.985.13"
Set-TextValue "E43" "  +0.66%  "

Set-TextValue "E44" "  -2.21%  "

Set-TextValue "E45" "  +0.62%  "

Set-TextValue "E46" "  +2.25%  "

Set-TextValue "D47" "9.25"
Set-TextValue "E47" "  +1.95%  "

Set-TextValue "D48" "2.712.53"
Set-TextValue "E48" "  -0.22%  "

Set-TextValue "D49" "97.78"
Set-TextValue "E49" "  -0.14%  "

Set-TextValue "D50" "67.62"
Set-TextValue "E50" "  -0.68%  "

Set-TextValue "D51" "73.86"
Set-TextValue "E51" "  +2.07%  "
